{"js": "// Replace each three-digit-by-one-digit multiplication expression in the\n// worksheet table with its updated counterpart. Every \"before\" text is\n// unique in the document, so searching the whole body for each literal\n// string and rewriting the single match it returns is sufficient.\n\nconst replacements = [\n  [\"378\u00d77=\", \"896\u00d72=\"],\n  [\"401\u00d77=\", \"723\u00d77=\"],\n  [\"468\u00d78=\", \"753\u00d77=\"],\n  [\"110\u00d76=\", \"318\u00d74=\"],\n  [\"414\u00d76=\", \"802\u00d75=\"],\n  [\"689\u00d79=\", \"183\u00d72=\"],\n  [\"662\u00d75=\", \"245\u00d79=\"],\n  [\"827\u00d75=\", \"723\u00d78=\"],\n  [\"972\u00d72=\", \"696\u00d78=\"],\n  [\"152\u00d78=\", \"504\u00d72=\"],\n  [\"661\u00d75=\", \"714\u00d78=\"],\n  [\"255\u00d72=\", \"728\u00d76=\"],\n  [\"547\u00d79=\", \"285\u00d76=\"],\n  [\"569\u00d79=\", \"904\u00d72=\"],\n  [\"653\u00d75=\", \"925\u00d78=\"],\n  [\"441\u00d73=\", \"422\u00d79=\"],\n  [\"533\u00d72=\", \"988\u00d72=\"],\n  [\"450\u00d72=\", \"102\u00d73=\"],\n  [\"225\u00d73=\", \"325\u00d78=\"],\n  [\"516\u00d78=\", \"184\u00d79=\"],\n  [\"534\u00d76=\", \"405\u00d76=\"],\n  [\"606\u00d77=\", \"358\u00d72=\"],\n  [\"210\u00d77=\", \"685\u00d78=\"],\n  [\"991\u00d77=\", \"505\u00d78=\"],\n  [\"219\u00d77=\", \"878\u00d74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication expression in the\n# worksheet table with its updated counterpart. Every occurrence of each\n# \"before\" text is unique in the document, so a set of literal\n# Find/Replace passes over the whole document body is sufficient and safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"378\u00d77=\"; New = \"896\u00d72=\" },\n    @{ Old = \"401\u00d77=\"; New = \"723\u00d77=\" },\n    @{ Old = \"468\u00d78=\"; New = \"753\u00d77=\" },\n    @{ Old = \"110\u00d76=\"; New = \"318\u00d74=\" },\n    @{ Old = \"414\u00d76=\"; New = \"802\u00d75=\" },\n    @{ Old = \"689\u00d79=\"; New = \"183\u00d72=\" },\n    @{ Old = \"662\u00d75=\"; New = \"245\u00d79=\" },\n    @{ Old = \"827\u00d75=\"; New = \"723\u00d78=\" },\n    @{ Old = \"972\u00d72=\"; New = \"696\u00d78=\" },\n    @{ Old = \"152\u00d78=\"; New = \"504\u00d72=\" },\n    @{ Old = \"661\u00d75=\"; New = \"714\u00d78=\" },\n    @{ Old = \"255\u00d72=\"; New = \"728\u00d76=\" },\n    @{ Old = \"547\u00d79=\"; New = \"285\u00d76=\" },\n    @{ Old = \"569\u00d79=\"; New = \"904\u00d72=\" },\n    @{ Old = \"653\u00d75=\"; New = \"925\u00d78=\" },\n    @{ Old = \"441\u00d73=\"; New = \"422\u00d79=\" },\n    @{ Old = \"533\u00d72=\"; New = \"988\u00d72=\" },\n    @{ Old = \"450\u00d72=\"; New = \"102\u00d73=\" },\n    @{ Old = \"225\u00d73=\"; New = \"325\u00d78=\" },\n    @{ Old = \"516\u00d78=\"; New = \"184\u00d79=\" },\n    @{ Old = \"534\u00d76=\"; New = \"405\u00d76=\" },\n    @{ Old = \"606\u00d77=\"; New = \"358\u00d72=\" },\n    @{ Old = \"210\u00d77=\"; New = \"685\u00d78=\" },\n    @{ Old = \"991\u00d77=\"; New = \"505\u00d78=\" },\n    @{ Old = \"219\u00d77=\"; New = \"878\u00d74=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute([ref]$r.Old, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$r.New, 2)\n}\n"}
